# Generate Report for Handoff
#
# Replaces the stale ".md" GUID-named entry with a fresh one and appends two
# new ".png" dependency rows (one per localization sheet), refreshing the
# handoff/handback timestamps to match the new run.

$wb = $excel.ActiveWorkbook

$newGuid = "334754f1-656a-41c1-887d-fe06454e87f0"
$newHash = "b4ef5fce6822c9e22b4465e98188a4be27cf779f"

$pngGuid1 = "5f78aa3f-3228-42d6-a1a2-09c67b271cb9"
$pngGuid2 = "fce06892-0e32-4d00-beb6-67c0ba47d31f"

$pngHash1 = "f4c1c50e61e9257c190a0eb9c7e91813f2ded374"
$pngHash2 = "96fa94eb92df8748ed65680744649ced076e638a"

$readyStatus = "Ready for handoff"
$overviewDate = "2016-51-12 20:51:13"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/175e83918ac4b1df8e6809b73a0f8ae720b44adf/e2e/"

# ---------------------------------------------------------------------------
# Sheet 1 : Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()

$overviewRows = @(
    @{ Row = 2; Name = ($newGuid + ".md") },
    @{ Row = 3; Name = ($pngGuid1 + ".png") },
    @{ Row = 4; Name = ($pngGuid2 + ".png") }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $ws1.Range("A" + $row).Value = $r.Name
    $ws1.Hyperlinks.Add($ws1.Range("A" + $row), ($mdBase + $r.Name), "", "", $r.Name)
    $ws1.Range("B" + $row).Value = $readyStatus
    $ws1.Range("C" + $row).Value = $readyStatus
    $ws1.Range("D" + $row).Value = $overviewDate
}

# ---------------------------------------------------------------------------
# Sheets 2 & 3 : per-locale detail (zh-cn / de-de)
# ---------------------------------------------------------------------------
$locales = @(
    @{ SheetIndex = 2; Locale = "zh-cn"; XlfHash = $newHash; Date = "2016-03-12 20:51:10";
       HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc894bff0ed8c4c913d0a41732b65b1d6d59cad8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" },
    @{ SheetIndex = 3; Locale = "de-de"; XlfHash = $newHash; Date = "2016-03-12 20:51:13";
       HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dba35b8bcebc9e08952d0fbfb4b4db27f87f010a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.SheetIndex)
    $ws.Hyperlinks.Delete()

    $mdName = $newGuid + ".md"
    $xlfName = $newGuid + "." + $loc.XlfHash + "." + $loc.Locale + ".xlf"
    $png1Name = $pngGuid1 + ".png"
    $png1TargetName = $pngHash1 + ".png"
    $png2Name = $pngGuid2 + ".png"
    $png2TargetName = $pngHash2 + ".png"

    # --- Row 2 : the refreshed .md entry -----------------------------------
    $ws.Range("A2").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("A2"), ($mdBase + $mdName), "", "", $mdName)

    $ws.Range("B2").Value = ".md"
    $ws.Hyperlinks.Add($ws.Range("B2"), ($mdBase + $mdName), "", "", ".md")

    $ws.Range("C2").Value = $readyStatus

    $ws.Range("D2").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("D2"), ($loc.HandoffBase + $xlfName), "", "", $xlfName)

    $ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("E2").Value = $loc.Date
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "Include"

    # --- Row 3 : first new .png dependency ---------------------------------
    $ws.Range("A3").Value = $png1Name
    $ws.Hyperlinks.Add($ws.Range("A3"), ($mdBase + $png1Name), "", "", $png1Name)

    $ws.Range("B3").Value = ".png"
    $ws.Hyperlinks.Add($ws.Range("B3"), ($mdBase + $png1Name), "", "", ".png")

    $ws.Range("C3").Value = $readyStatus

    $ws.Range("D3").Value = $png1TargetName
    $ws.Hyperlinks.Add($ws.Range("D3"), ($loc.HandoffBase + $png1TargetName), "", "", $png1TargetName)

    $ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("E3").Value = $loc.Date
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "IsDependency"
    $ws.Range("J3").Value = "e2e\" + $mdName

    # --- Row 4 : second new .png dependency --------------------------------
    $ws.Range("A4").Value = $png2Name
    $ws.Hyperlinks.Add($ws.Range("A4"), ($mdBase + $png2Name), "", "", $png2Name)

    $ws.Range("B4").Value = ".png"
    $ws.Hyperlinks.Add($ws.Range("B4"), ($mdBase + $png2Name), "", "", ".png")

    $ws.Range("C4").Value = $readyStatus

    $ws.Range("D4").Value = $png2TargetName
    $ws.Hyperlinks.Add($ws.Range("D4"), ($loc.HandoffBase + $png2TargetName), "", "", $png2TargetName)

    $ws.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("E4").Value = $loc.Date
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "IsDependency"
    $ws.Range("J4").Value = "e2e\" + $mdName
}
